$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: points (F) 21 -> 22, matches (G) 15 -> 16
$ws.Range("F11").Value = 22
$ws.Range("G11").Value = 16

# Row 13: position (E) 12 -> 13
$ws.Range("E13").Value = 13

# Row 14: position (E) 13 -> 12, points (F) 19 -> 20, matches (G) 14 -> 15
$ws.Range("E14").Value = 12
$ws.Range("F14").Value = 20
$ws.Range("G14").Value = 15
